# Auto-generated Excel COM-interop script
# Applies the cached numeric-value refresh described in the commit diff
# ("chore: update Sheets via scheduled runner") across all 8 Leve-profit sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 368
$ws.Range("I2").Value = 310.54544
$ws.Range("K2").Value = 310.54544
$ws.Range("M2").Value = -197.54544
$ws.Range("H17").Value = 2292367
$ws.Range("J17").Value = 2292367
$ws.Range("L17").Value = 6877101
$ws.Range("N17").Value = -6877437
$ws.Range("H70").Value = 2302.8333
$ws.Range("I70").Value = 2362.75
$ws.Range("J70").Value = 2254.9
$ws.Range("K70").Value = 7088.25
$ws.Range("L70").Value = 6764.700000000001
$ws.Range("M70").Value = -6818.25
$ws.Range("N70").Value = -7304.700000000001
$ws.Range("H73").Value = 2302.8333
$ws.Range("I73").Value = 2362.75
$ws.Range("J73").Value = 2254.9
$ws.Range("K73").Value = 7088.25
$ws.Range("L73").Value = 6764.700000000001
$ws.Range("M73").Value = -6152.25
$ws.Range("N73").Value = -8636.700000000001
$ws.Range("H97").Value = 4975
$ws.Range("J97").Value = 4975
$ws.Range("L97").Value = 14925
$ws.Range("N97").Value = -15917
$ws.Range("H123").Value = 18485.625
$ws.Range("J123").Value = 18485.625
$ws.Range("L123").Value = 18485.625
$ws.Range("N123").Value = -28285.625
$ws.Range("H125").Value = 2694.2173
$ws.Range("I125").Value = 2987.4443
$ws.Range("J125").Value = 2505.7144
$ws.Range("K125").Value = 26886.9987
$ws.Range("L125").Value = 22551.4296
$ws.Range("M125").Value = -24426.9987
$ws.Range("N125").Value = -27471.4296
$ws.Range("H129").Value = 1100.4032
$ws.Range("J129").Value = 1107.2543
$ws.Range("L129").Value = 3321.7629
$ws.Range("N129").Value = -13321.7629
$ws.Range("H132").Value = 2347.2727
$ws.Range("I132").Value = 2165.077
$ws.Range("J132").Value = 2610.4443
$ws.Range("K132").Value = 6495.231000000001
$ws.Range("L132").Value = 7831.3329
$ws.Range("M132").Value = -3965.231000000001
$ws.Range("N132").Value = -12891.3329
$ws.Range("H135").Value = 2058.6135
$ws.Range("I135").Value = 1509.7273
$ws.Range("K135").Value = 13587.5457
$ws.Range("M135").Value = -11052.5457
$ws.Range("H138").Value = 2247.3735
$ws.Range("I138").Value = 964.12195
$ws.Range("K138").Value = 2892.36585
$ws.Range("M138").Value = 2247.63415

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5036.32
$ws.Range("I32").Value = 3390.1013
$ws.Range("J32").Value = 11229.238
$ws.Range("K32").Value = 3390.1013
$ws.Range("L32").Value = 11229.238
$ws.Range("M32").Value = -3103.1013
$ws.Range("N32").Value = -11803.238
$ws.Range("H45").Value = 11108.4
$ws.Range("I45").Value = 11108.4
$ws.Range("K45").Value = 11108.4
$ws.Range("M45").Value = -10731.4
$ws.Range("H61").Value = 4850.625
$ws.Range("I61").Value = 6622.4
$ws.Range("K61").Value = 6622.4
$ws.Range("M61").Value = -6410.4
$ws.Range("H63").Value = 1000000000
$ws.Range("I63").Value = 1000000000
$ws.Range("K63").Value = 1000000000
$ws.Range("M63").Value = -999999314
$ws.Range("H66").Value = 1000000000
$ws.Range("I66").Value = 1000000000
$ws.Range("K66").Value = 5000000000
$ws.Range("M66").Value = -4999996568
$ws.Range("H74").Value = 966.44116
$ws.Range("I74").Value = 814.3674
$ws.Range("J74").Value = 1358.6316
$ws.Range("K74").Value = 814.3674
$ws.Range("L74").Value = 1358.6316
$ws.Range("M74").Value = 59.63260000000002
$ws.Range("N74").Value = -3106.6316
$ws.Range("H77").Value = 966.44116
$ws.Range("I77").Value = 814.3674
$ws.Range("J77").Value = 1358.6316
$ws.Range("K77").Value = 4071.837
$ws.Range("L77").Value = 6793.157999999999
$ws.Range("M77").Value = 296.163
$ws.Range("N77").Value = -15529.158
$ws.Range("H110").Value = 710
$ws.Range("I110").Value = 636.25
$ws.Range("J110").Value = 1300
$ws.Range("K110").Value = 636.25
$ws.Range("L110").Value = 1300
$ws.Range("M110").Value = 1408.75
$ws.Range("N110").Value = -5390
$ws.Range("H132").Value = 2624.76
$ws.Range("I132").Value = 1741.8108
$ws.Range("J132").Value = 5137.769
$ws.Range("K132").Value = 5225.4324
$ws.Range("L132").Value = 15413.307
$ws.Range("M132").Value = -2695.4324
$ws.Range("N132").Value = -20473.307
$ws.Range("H136").Value = 4850.625
$ws.Range("I136").Value = 6622.4
$ws.Range("K136").Value = 19867.2
$ws.Range("M136").Value = -17317.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1689
$ws.Range("I86").Value = 1645.5385
$ws.Range("J86").Value = 1783.1666
$ws.Range("K86").Value = 1645.5385
$ws.Range("L86").Value = 1783.1666
$ws.Range("M86").Value = -522.5385000000001
$ws.Range("N86").Value = -4029.1666
$ws.Range("H89").Value = 1689
$ws.Range("I89").Value = 1645.5385
$ws.Range("J89").Value = 1783.1666
$ws.Range("K89").Value = 8227.692500000001
$ws.Range("L89").Value = 8915.833000000001
$ws.Range("M89").Value = -2611.692500000001
$ws.Range("N89").Value = -20147.833
$ws.Range("H134").Value = 4940.6216
$ws.Range("I134").Value = 6818.1
$ws.Range("K134").Value = 20454.3
$ws.Range("M134").Value = -17919.3

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1800.7273
$ws.Range("I12").Value = 1401.1428
$ws.Range("J12").Value = 2500
$ws.Range("K12").Value = 1401.1428
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = -1231.1428
$ws.Range("N12").Value = -2840
$ws.Range("H132").Value = 2424.2827
$ws.Range("I132").Value = 2066.3667
$ws.Range("J132").Value = 3095.375
$ws.Range("K132").Value = 6199.1001
$ws.Range("L132").Value = 9286.125
$ws.Range("M132").Value = -3669.1001
$ws.Range("N132").Value = -14346.125
$ws.Range("H134").Value = 1915.862
$ws.Range("I134").Value = 2200.5366
$ws.Range("J134").Value = 1229.2941
$ws.Range("K134").Value = 6601.6098
$ws.Range("L134").Value = 3687.8823
$ws.Range("M134").Value = -4066.6098
$ws.Range("N134").Value = -8757.882300000001
$ws.Range("H141").Value = 38876
$ws.Range("J141").Value = 38876
$ws.Range("L141").Value = 38876
$ws.Range("N141").Value = -49236

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 209188
$ws.Range("I5").Value = 314.56
$ws.Range("J5").Value = 436224.34
$ws.Range("K5").Value = 943.6800000000001
$ws.Range("L5").Value = 1308673.02
$ws.Range("M5").Value = -831.6800000000001
$ws.Range("N5").Value = -1308897.02
$ws.Range("H122").Value = 3353.7354
$ws.Range("I122").Value = 442.2
$ws.Range("J122").Value = 4566.875
$ws.Range("K122").Value = 3979.8
$ws.Range("L122").Value = 41101.875
$ws.Range("M122").Value = -1529.8
$ws.Range("N122").Value = -46001.875
$ws.Range("H135").Value = 209188
$ws.Range("I135").Value = 314.56
$ws.Range("J135").Value = 436224.34
$ws.Range("K135").Value = 2831.04
$ws.Range("L135").Value = 3926019.06
$ws.Range("M135").Value = -296.04
$ws.Range("N135").Value = -3931089.06

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1807.2742
$ws.Range("I132").Value = 1379.4615
$ws.Range("K132").Value = 4138.3845
$ws.Range("M132").Value = -1608.3845

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 11136224
$ws.Range("I132").Value = 14445890
$ws.Range("J132").Value = 3709.3635
$ws.Range("K132").Value = 43337670
$ws.Range("L132").Value = 11128.0905
$ws.Range("M132").Value = -43335140
$ws.Range("N132").Value = -16188.0905

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 1000000000
$ws.Range("I12").Value = 1000000000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1000000000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -999999858
$ws.Range("N12").ClearContents()
$ws.Range("H123").Value = 28971.428
$ws.Range("J123").Value = 28971.428
$ws.Range("L123").Value = 28971.428
$ws.Range("N123").Value = -38771.428
$ws.Range("H132").Value = 19548.5
$ws.Range("I132").Value = 23465.818
$ws.Range("J132").Value = 2312.3
$ws.Range("K132").Value = 70397.454
$ws.Range("L132").Value = 6936.900000000001
$ws.Range("M132").Value = -67867.454
$ws.Range("N132").Value = -11996.9
$ws.Range("H136").Value = 7044513
$ws.Range("I136").Value = 2283.5532
$ws.Range("J136").Value = 20835546
$ws.Range("K136").Value = 6850.659599999999
$ws.Range("L136").Value = 62506638
$ws.Range("M136").Value = -4300.659599999999
$ws.Range("N136").Value = -62511738

Write-Output "Applied 215 cell updates across 8 sheets."